$wb = $excel.ActiveWorkbook

# Update match score predictions for Victor, Gabrielzinho and Matheus

$ws = $wb.Worksheets.Item("Victor")
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 6).Value = 2
$ws.Cells.Item(16, 4).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 6).Value = 2
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 6).Value = 3
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 6).Value = 2
$ws.Cells.Item(33, 4).Value = 2
$ws.Cells.Item(33, 6).Value = 2
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(40, 6).Value = 3
$ws.Cells.Item(41, 4).Value = 2
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(48, 4).Value = 2
$ws.Cells.Item(48, 6).Value = 1
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 6).Value = 2
$ws.Cells.Item(56, 4).Value = 3
$ws.Cells.Item(56, 6).Value = 3
$ws.Cells.Item(57, 4).Value = 1
$ws.Cells.Item(57, 6).Value = 2
$ws.Cells.Item(64, 4).Value = 2
$ws.Cells.Item(64, 6).Value = 1
$ws.Cells.Item(65, 4).Value = 2
$ws.Cells.Item(65, 6).Value = 3

$ws = $wb.Worksheets.Item("Matheus")
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 6).Value = 2
$ws.Cells.Item(16, 4).Value = 2
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 6).Value = 2
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 6).Value = 2
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 6).Value = 2
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 6).Value = 2
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(48, 6).Value = 1
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 6).Value = 3
$ws.Cells.Item(56, 4).Value = 1
$ws.Cells.Item(56, 6).Value = 1
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 6).Value = 2
$ws.Cells.Item(64, 4).Value = 2
$ws.Cells.Item(64, 6).Value = 1
$ws.Cells.Item(65, 4).Value = 1
$ws.Cells.Item(65, 6).Value = 2

$ws = $wb.Worksheets.Item("Gabrielzinho")
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(16, 4).Value = 2
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 6).Value = 2
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 6).Value = 2
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 6).Value = 2
$ws.Cells.Item(33, 4).Value = 2
$ws.Cells.Item(33, 6).Value = 2
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 6).Value = 2
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(48, 4).Value = 2
$ws.Cells.Item(48, 6).Value = 1
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 6).Value = 3
$ws.Cells.Item(56, 4).Value = 1
$ws.Cells.Item(56, 6).Value = 3
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 6).Value = 1
$ws.Cells.Item(64, 4).Value = 2
$ws.Cells.Item(64, 6).Value = 1
$ws.Cells.Item(65, 4).Value = 2
$ws.Cells.Item(65, 6).Value = 3

# Move to the sheet that was last edited, matching the workbooks active-tab state
$gab = $wb.Worksheets.Item("Gabrielzinho")
$gab.Activate()
$gab.Range("D66").Select()
